# Level & Industry integrated
# Adds 4 new rows (259-262) to both the "en" and "de" sheets with
# new localization keys/values for Industry / HierarchyLevel terms.

$wb = $excel.ActiveWorkbook

$enSheet = $wb.Worksheets.Item("en")
$deSheet = $wb.Worksheets.Item("de")

function Set-LocRow {
    param(
        $ws,
        [int]$row,
        [string]$keyValue,
        [string]$localizedValue
    )

    $keyCell = $ws.Cells.Item($row, 1)
    $keyCell.Value = $keyValue
    $keyCell.WrapText = $true
    $keyCell.VerticalAlignment = -4108   # xlCenter

    $valCell = $ws.Cells.Item($row, 2)
    $valCell.Value = $localizedValue
    $valCell.WrapText = $true
    $valCell.NumberFormat = "0.00"
}

# "en" sheet (column A and B both in English)
Set-LocRow $enSheet 259 "Industry"       "Industry"
Set-LocRow $enSheet 260 "HierarchyLevel" "Hierarchy level"
Set-LocRow $enSheet 261 "Industries"     "Industries"
Set-LocRow $enSheet 262 "HierarchyLevels" "Hierarchy levels"

# "de" sheet (column A English key, column B German translation)
Set-LocRow $deSheet 259 "Industry"        "Industrie"
Set-LocRow $deSheet 260 "Level"           "Hierarchiestufe"
Set-LocRow $deSheet 261 "Industries"      "Industrien"
Set-LocRow $deSheet 262 "HierarchyLevels" "Hierarchiestufen"

# Column width adjustments observed in the diff
$enSheet.Columns.Item(1).ColumnWidth = 34.333333333333336
$deSheet.Columns.Item(1).ColumnWidth = 34.333333333333336

# Update selection / view state to match the authored edit
[void]$enSheet.Range("A261:B262").Select()

[void]$deSheet.Range("A261:B262").Select()
